$wb = $excel.ActiveWorkbook

# Each sheet contains a block of leading data rows (rows 2..N) that are
# exact duplicates of the last rows of the same sheet. The edit removes
# those leading duplicate rows; every row below shifts up accordingly,
# which also shrinks each sheet's used range / <dimension>.

$ws1 = $wb.Worksheets.Item("s__Ligilactobacillus animalis-b-p")
$ws1.Range("A2:A4").EntireRow.Delete()

$ws2 = $wb.Worksheets.Item("s__Ligilactobacillus murinus-b-p")
$ws2.Range("A2:A17").EntireRow.Delete()

$ws3 = $wb.Worksheets.Item("s__Ligilactobacillus ruminis-b-p")
$ws3.Range("A2:A2").EntireRow.Delete()

$ws4 = $wb.Worksheets.Item("s__Ligilactobacillus salivarius-b-p")
$ws4.Range("A2:A5").EntireRow.Delete()
